$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 18 de Mayo de 2020 a las 01:05'
$ws.Range("B4").Value = 1526807
$ws.Range("C4").Value = 19034
$ws.Range("D4").Value = 344933
$ws.Range("E4").Value = 1090901
$ws.Range("G4").Value = 860
$ws.Range("H4").Value = 90973
$ws.Range("B8").Value = 241080
$ws.Range("C8").Value = 7938
$ws.Range("D8").Value = 94122
$ws.Range("E8").Value = 130840
$ws.Range("G8").Value = 485
$ws.Range("H8").Value = 16118
$ws.Range("A39").Value = 'Japon'
$ws.Range("B39").Value = 16285
$ws.Range("C39").Value = 48
$ws.Range("D39").Value = 11153
$ws.Range("E39").Value = 4388
$ws.Range("G39").Value = 19
$ws.Range("H39").Value = 744
$ws.Range("A40").Value = 'Austria'
$ws.Range("B40").Value = 16242
$ws.Range("C40").Value = 41
$ws.Range("D40").Value = 14563
$ws.Range("E40").Value = 1050
$ws.Range("H40").Value = 629
$ws.Range("A41").Value = 'Colombia'
$ws.Range("B41").Value = 15574
$ws.Range("C41").Value = 635
$ws.Range("D41").Value = 3751
$ws.Range("E41").Value = 11249
$ws.Range("G41").Value = 12
$ws.Range("H41").Value = 574
$ws.Range("A42").Value = 'Sudafrica'
$ws.Range("B42").Value = 15515
$ws.Range("C42").Value = 1160
$ws.Range("D42").Value = 7006
$ws.Range("E42").Value = 8245
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 264
$ws.Range("B52").Value = 8249
$ws.Range("C52").Value = 12
$ws.Range("E52").Value = 7985
$ws.Range("A63").Value = 'Nigeria'
$ws.Range("B63").Value = 5959
$ws.Range("C63").Value = 338
$ws.Range("D63").Value = 1594
$ws.Range("E63").Value = 4183
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 182
$ws.Range("A64").Value = 'Ghana'
$ws.Range("B64").Value = 5735
$ws.Range("D64").Value = 1754
$ws.Range("E64").Value = 3952
$ws.Range("H64").Value = 29
$ws.Range("A125").Value = 'Venezuela'
$ws.Range("B125").Value = 541
$ws.Range("C125").Value = 37
$ws.Range("D125").Value = 241
$ws.Range("E125").Value = 290
$ws.Range("H125").Value = 10
$ws.Range("A126").Value = 'Jamaica'
$ws.Range("B126").Value = 517
$ws.Range("C126").Value = 6
$ws.Range("D126").Value = 121
$ws.Range("E126").Value = 387
$ws.Range("H126").Value = 9
$ws.Range("A127").Value = 'Tanzania'
$ws.Range("B127").Value = 509
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 183
$ws.Range("E127").Value = 305
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 21
$ws.Range("A128").Value = 'Sierra Leona'
$ws.Range("B128").Value = 505
$ws.Range("C128").Value = 43
$ws.Range("D128").Value = 141
$ws.Range("E128").Value = 332
$ws.Range("G128").Value = 3
$ws.Range("H128").Value = 32
$ws.Range("A130").Value = 'Haiti'
$ws.Range("B130").Value = 456
$ws.Range("C130").Value = 98
$ws.Range("D130").Value = 21
$ws.Range("E130").Value = 415
$ws.Range("H130").Value = 20
$ws.Range("A131").Value = 'Reunion'
$ws.Range("B131").Value = 443
$ws.Range("D131").Value = 354
$ws.Range("E131").Value = 89
$ws.Range("H131").Value = 0
$ws.Range("A132").Value = 'Taiwan'
$ws.Range("B132").Value = 440
$ws.Range("D132").Value = 395
$ws.Range("E132").Value = 38
$ws.Range("H132").Value = 7
$ws.Range("A133").Value = 'Congo'
$ws.Range("B133").Value = 391
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 87
$ws.Range("E133").Value = 289
$ws.Range("H133").Value = 15
$ws.Range("A134").Value = 'Estado de Palestina'
$ws.Range("B134").Value = 381
$ws.Range("C134").Value = 5
$ws.Range("D134").Value = 335
$ws.Range("E134").Value = 44
$ws.Range("H134").Value = 2
